# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.
#
# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"):
#   - fix the E1 header label (it mistakenly held a stray numeric value
#     instead of a "2050"/"2041-2050" text label, like the other year
#     headers in B1:D1)
#   - remove the "Total" row (row 13) at the bottom of the table
# Sheet 5 ("Emissoes Totais (MtCO2eq)") only needs the E1 header label fixed
#   (it never had a Total row).
# Sheet 6 ("Custo Total (bilhoes de R$)") only needs its "Total" row (row 4)
#   removed (its header row has no year labels).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Potencia Acumulada - SIN (MW) ---
$ws1 = $wb.Worksheets.Item(1)
# Format as text first so the numeric-looking "2050" is stored as a literal
# string rather than being re-interpreted as a number.
$ws1.Range("E1").NumberFormat = "@"
$ws1.Range("E1").Value = "2050"
$ws1.Range("A13:E13").EntireRow.Delete()

# --- Sheet 2: Geracao Periodo Medio (MWMed) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E1").NumberFormat = "@"
$ws2.Range("E1").Value = "2050"
$ws2.Range("A13:E13").EntireRow.Delete()

# --- Sheet 3: Atendimento a Ponta(MW) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E1").NumberFormat = "@"
$ws3.Range("E1").Value = "2050"
$ws3.Range("A13:E13").EntireRow.Delete()

# --- Sheet 4: Potencia Incremental - SIN(MW) ---
# This sheet's header row uses ranges (2015-2030, 2031-2040), so the fixed
# label follows the same pattern: 2041-2050 (already non-numeric-looking,
# so it is stored as text automatically).
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "2041-2050"
$ws4.Range("A13:E13").EntireRow.Delete()

# --- Sheet 5: Emissoes Totais (MtCO2eq) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").NumberFormat = "@"
$ws5.Range("E1").Value = "2050"

# --- Sheet 6: Custo Total (bilhoes de R$) ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A4:B4").EntireRow.Delete()
